$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Tasks done / hours for the existing row 12 (2024-12-02)
$ws.Range("B12").Value = "fixing some bugs in the combat system and implementing it into the main area, some more visual for the inventory and the walking animations"
$ws.Range("C12").Value = 6

# Add a new row 13 for 2024-12-03 (serial date number 45629)
$ws.Range("A13").Value = 45629
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
$ws.Range("B13").Value = "implemented money and most of the systems required so far for the hot/cold puzzle types"

# Update selection to match the final state
$ws.Range("B13").Select()
